# "Add files via upload" — append a new response row (row 4) to the
# SCAP submissions sheet, matching the data of an existing row's
# formatting (styles/number formats) but with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 already has the correct per-column formatting (text style,
# datetime style, plain-number style, date style) for every column
# A:V, so clone it into row 4 first and then overwrite the values.
$ws.Range("A3:V3").Copy()
$ws.Range("A4:V4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A4").Value = "Pba2pP"
$ws.Range("B4").Value = "jBVv5Q"
$ws.Range("C4").Value = 45436.802708333336
$ws.Range("D4").Value = "andre.amorim@planejamento.mg.gov.br"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "André"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 45422
$ws.Range("K4").Value = 45429
$ws.Range("L4").Value = "EPPGG"
$ws.Range("M4").Value = "III"
$ws.Range("N4").Value = "EPPGG"
$ws.Range("O4").Value = "C"
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 45436
$ws.Range("S4").Value = 46897
$ws.Range("T4").Value = 45431
$ws.Range("U4").Value = "V"
$ws.Range("V4").Value = "G"
